# Presentacion_juegos_en_red_fase_4.pptx - "Arreglo de las presentaciones"
#
# - La presentacion de la fase 4 se deja exclusiva de esta fase (se fusionan
#   los runs "FASE "/"4" -> "FASE 4" en la caja de texto de la portada y en
#   el titulo "AÑADIDOS FASE 4").
# - Se sustituye el contenido de la diapositiva de añadidos por los puntos
#   reales de la fase 4 (lobby sin IP + websockets).

$p = $ppt.ActivePresentation

# --- Slide 1: portada "FASE 4" -------------------------------------------
$s1 = $p.Slides.Item(1)
$fase4Box = $s1.Shapes.Item(8)   # "CuadroTexto 2" -> runs "FASE " + "4"
$tr = $fase4Box.TextFrame.TextRange
# Forzar un cambio real de caracteres para que el motor funda los dos runs
# "FASE " / "4" en un unico run con el texto final "FASE 4".
$tr.Text = "FASE 4_"
$tr.Text = "FASE 4"

# --- Slide 2: "AÑADIDOS FASE 4" -------------------------------------------
$s2 = $p.Slides.Item(2)

$title2 = $s2.Shapes.Item(1)     # "Titulo 1" -> runs "AÑADIDOS FASE " + "4"
$trTitle2 = $title2.TextFrame.TextRange
$trTitle2.Text = "AÑADIDOS FASE 4_"
$trTitle2.Text = "AÑADIDOS FASE 4"

$body2 = $s2.Shapes.Item(2)      # "Marcador de contenido 2" -> lista de viñetas
$trBody2 = $body2.TextFrame.TextRange
$trBody2.Text = "ENTRADA AL LOBBY SIN INTRODUCIR IP, Y USO DEL MISMO PARA INICIAR UNA PARTIDA.`rUSO DE WEBSOCKETS PARA JUGAR UNA PARTIDA ENTRE 2 JUGADORES DESDE DISTINTOS ORDENADORES.`rUSO DE WEBSOCKETS PARA EL INTERCAMBIO DEL ESTADO DE LA PARTIDA ENTRE LOS JUGADORES."
